$d = $word.ActiveDocument

$d.Content.Find.Execute("39÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷9=", 2) | Out-Null
$d.Content.Find.Execute("80÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷8=", 2) | Out-Null
$d.Content.Find.Execute("67÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷7=", 2) | Out-Null
$d.Content.Find.Execute("95÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷4=", 2) | Out-Null
$d.Content.Find.Execute("13÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷7=", 2) | Out-Null
$d.Content.Find.Execute("56÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷4=", 2) | Out-Null
$d.Content.Find.Execute("58÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷4=", 2) | Out-Null
$d.Content.Find.Execute("10÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷4=", 2) | Out-Null
$d.Content.Find.Execute("15÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷3=", 2) | Out-Null
$d.Content.Find.Execute("95÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷3=", 2) | Out-Null
$d.Content.Find.Execute("39÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷6=", 2) | Out-Null
$d.Content.Find.Execute("40÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷4=", 2) | Out-Null
$d.Content.Find.Execute("46÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷6=", 2) | Out-Null
$d.Content.Find.Execute("55÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷8=", 2) | Out-Null
$d.Content.Find.Execute("58÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷9=", 2) | Out-Null
$d.Content.Find.Execute("69÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷6=", 2) | Out-Null
$d.Content.Find.Execute("77÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷4=", 2) | Out-Null
$d.Content.Find.Execute("52÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷7=", 2) | Out-Null
$d.Content.Find.Execute("53÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷8=", 2) | Out-Null
$d.Content.Find.Execute("14÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷5=", 2) | Out-Null
$d.Content.Find.Execute("59÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷8=", 2) | Out-Null
$d.Content.Find.Execute("33÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "43÷4=", 2) | Out-Null
$d.Content.Find.Execute("44÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷7=", 2) | Out-Null
$d.Content.Find.Execute("16÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷5=", 2) | Out-Null
$d.Content.Find.Execute("47÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷7=", 2) | Out-Null

"Done"